$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Efna1","Epha5","ECs",2,1,18.819997,37.639994,0.922622403444,0.8996825797928821,1,0.5,0.006466,0.012932,0.008493427970384656,0.008493427970384656,0.121690100602,0.4867604024080001,0.007836226927514787,0.00764138918768069),
    @("ECs","Efna1","Epha5","MuSCs",2,1,18.819997,37.639994,0.922622403444,0.8996825797928821,2,1,0.7548284999999999,1.509657,0.9915065720296153,0.9915065720296153,14.2058701055145,56.823480422058,0.9147861765164852,0.8920411906052015),
    @("FAPs","Efna1","Epha5","ECs",3,1,1.040222666666667,3.120668,0.05099537140398731,0.07459115527268931,1,0.5,0.006466,0.012932,0.008493427970384656,0.008493427970384656,0.006726079762666668,0.04035647857600001,0.0004331255138427796,0.0006335346045363643),
    @("FAPs","Efna1","Epha5","MuSCs",3,1,1.040222666666667,3.120668,0.05099537140398731,0.07459115527268931,2,1,0.7548284999999999,1.509657,0.9915065720296153,0.9915065720296153,0.785189715146,4.711138290876,0.05056224589014453,0.07395762066815294),
    @("MuSCs","Efna1","Epha5","ECs",2,1,0.5381545,1.076309,0.02638222515201273,0.02572626493442845,1,0.5,0.006466,0.012932,0.008493427970384656,0.008493427970384656,0.003479706997,0.013918827988,0.0002240755290270905,0.0002185041781676006),
    @("MuSCs","Efna1","Epha5","MuSCs",2,1,0.5381545,1.076309,0.02638222515201273,0.02572626493442845,2,1,0.7548284999999999,1.509657,0.9915065720296153,0.9915065720296153,0.40621435400325,1.624857416013,0.02615814962298564,0.02550776075626085)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
